# Penalty Reward System (unfinished) - remove some rows from the
# "Weekly Quantity" and "Monthly Trend" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Quantity": delete rows 20 and 21 (weeks of 2024-03-10
# and 2024-03-17), shifting the remaining rows up. ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows.Item(20).Delete()
$wsWeekly.Rows.Item(20).Delete()

# --- Sheet "Monthly Trend": delete row 8 (month ending 2024-03-31),
# shifting the remaining row up. ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Rows.Item(8).Delete()
